# 202309 September full release
#
# The "ANSP" sheet lists ANSP names alphabetically in column A. A new
# ANSP ("BHANSA") needs to be inserted between "Avinor (Continental)"
# (row 6) and "BULATSA" (row 7), pushing every row below it down by one.

$wb = $excel.ActiveWorkbook
$dataWs = $wb.Worksheets.Item("ACE_landing_page_data")
$ws = $wb.Worksheets.Item("ANSP")

# Shift the existing body rows (A7:A39) down by one row. Copying the
# range and inserting it back at A7 (shift down) keeps each row's own
# formatting/borders attached to its value as it moves - unlike a plain
# blank row insert, which only carries over the row above's font/border.
$ws.Range("A7:A39").Copy() | Out-Null
$ws.Range("A7").Insert(-4121) | Out-Null   # xlShiftDown

# Fill the now-empty A7 with the new entry.
$ws.Cells.Item(7, 1).Value2 = "BHANSA"

# Give it the same "middle of list" look as its neighbours: a thin grey
# left border (the other body rows use this single-sided border; only
# the very last row also gets a matching bottom border).
$leftBorder = $ws.Cells.Item(7, 1).Borders.Item(7)   # xlEdgeLeft
$leftBorder.LineStyle = 1                            # xlContinuous
$leftBorder.Weight = 2                               # xlThin
$leftBorder.Color = 10066329                         # RGB(153,153,153)

# Reflect where the user was last working on the ANSP sheet...
$ws.Range("A8").Select()

# ...but leave the landing-page data sheet as the active/visible tab,
# matching the workbook's saved state.
$dataWs.Activate()
